# Adds a new "Boolean column" / "Nullable Boolean column" pair of columns
# to both worksheets (configuration for boolean strings), renaming the
# existing "Bool column" header to "Boolean column", and re-selects the
# "Data OK" worksheet (and cell K3 on each sheet) to match the author's
# final view state.

$wb = $excel.ActiveWorkbook

$wsOk  = $wb.Worksheets.Item("Data OK")
$wsErr = $wb.Worksheets.Item("Data With Errors")

# ---------------------------------------------------------------------
# "Data OK" sheet (sheet1.xml)
# ---------------------------------------------------------------------

# Rename the existing boolean header and add the new nullable-boolean one.
# (Nullable Boolean column is written first so it claims the earlier shared
# string table slot, matching the author's original edit order.)
$wsOk.Range("K1").Value = "Nullable Boolean column"
$wsOk.Range("K1").Font.Bold = $true
$wsOk.Range("J1").Value = "Boolean column"

# Mirror the non-nullable boolean column values into the new column.
$wsOk.Range("K2").Value = 1
$wsOk.Range("K4").Value = 0
$wsOk.Range("K6").Value = "S"

# ---------------------------------------------------------------------
# "Data With Errors" sheet (sheet2.xml)
# ---------------------------------------------------------------------

$wsErr.Range("K1").Value = "Nullable Boolean column"
$wsErr.Range("K1").Font.Bold = $true
$wsErr.Range("J1").Value = "Boolean column"

$wsErr.Range("K2").Value = 1
$wsErr.Range("J3").Value = "SDDD"
$wsErr.Range("K3").Value = "DDD"
$wsErr.Range("K4").Value = 0
$wsErr.Range("K5").Value = "A"
$wsErr.Range("K6").Value = "S"

# The new column needed an explicit width/bestFit definition once it held
# data (mirrors the "Data OK" sheet's pre-existing column J metadata).
$wsErr.Columns.Item(10).ColumnWidth = 11.17

# ---------------------------------------------------------------------
# View state: "Data OK" becomes the active/selected tab with K3 selected,
# while "Data With Errors" keeps the same cursor position but loses the
# tab selection.
# ---------------------------------------------------------------------

[void]$wsErr.Range("K3").Select()
[void]$wsOk.Select()
[void]$wsOk.Range("K3").Select()
